$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text: "stage_Vx" -> "stage"
$ws.Range("G6").Value = "stage"

# Update column G (rows 7-30) - growth stage data changed from numeric 3 to text "R3"
for ($r = 7; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = "R3"
    $cell.NumberFormat = "0"
}

# Update the Notes cell (A4) with the additional note
$ws.Range("A4").Value = "Notes: Some phenology areas are in funky parts of C4 plots. Was wrong, corrected on Dec 17 2019"

# Update selection to A5
$ws.Range("A5").Select()
